# UnitTestDevelopmentPlan.xlsx — "added unit coices, unit tests and used factors for Volume"
#
# 1) Row 49 was a stray duplicate "Force Gradient" entry (A49 pointed at the
#    same shared string as A45) carrying a leftover Factors formula in P49.
#    Deleting the whole row shifts everything below it up by one, which is
#    exactly what the target file shows (A49 -> "Specific Heat Capacity
#    Temperature", the old P51 factor formula now lives in P50, the summary
#    row moves from 63 to 62, the "remains" row moves from 64 to 63, and the
#    trailing "Gamma Ray" row moves from 65 to 64).
# 2) The "Volume" unit-test row (row 46) gets Start/Finish test dates filled
#    in (B46/C46 = 8/13/2024, matching the same date used on the surrounding
#    rows), using the same date number format already used by column B/C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the duplicated "Force Gradient" row ---------------------------
$ws.Rows(49).Delete() | Out-Null

# --- Fill in the Volume unit-test Start/Finish dates -----------------------
$ws.Range("B44:C44").Copy() | Out-Null
$ws.Range("B46:C46").PasteSpecial(-4122) | Out-Null
$ws.Range("B46").Value = 45517
$ws.Range("C46").Value = 45517

# --- Restore the selection left by the editor ------------------------------
$ws.Range("D51").Select() | Out-Null
